$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nbsp = [char]0x00A0

# Row 4 = duplicate of row 2 (Sharjah / Kings XI Punjab match)
$ws.Range("A4").Value = " Sharjah"
$ws.Range("B4").Value = " October 26 2020"
$ws.Range("C4").Value = "Kings XI won by 8 wickets (with 7 balls remaining)"
$ws.Range("D4").Value = "Kolkata Knight Riders"
$ws.Range("E4").Value = "Kings XI Punjab"
$ws.Range("F4").Value = "Lockie Ferguson$nbsp"
$ws.Range("G4").Value = "'24"
$ws.Range("H4").Value = "'13"
$ws.Range("I4").Value = "'3"
$ws.Range("J4").Value = "'1"
$ws.Range("K4").Value = "'184.61"

# Row 5 = duplicate of row 3 (Abu Dhabi / RCB match)
$ws.Range("A5").Value = " Abu Dhabi"
$ws.Range("B5").Value = " October 21 2020"
$ws.Range("C5").Value = "RCB won by 8 wickets (with 39 balls remaining)"
$ws.Range("D5").Value = "Kolkata Knight Riders"
$ws.Range("E5").Value = "Royal Challengers Bangalore"
$ws.Range("F5").Value = "Lockie Ferguson$nbsp"
$ws.Range("G5").Value = "'19"
$ws.Range("H5").Value = "'16"
$ws.Range("I5").Value = "'1"
$ws.Range("J5").Value = "'0"
$ws.Range("K5").Value = "'118.75"
